# Add a 5th data row to sheet1, mirroring the style/format of existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain-text / number cells (no special style needed - matches B4/D4 pattern
# but row 5 keeps them unstyled per target).
$ws.Range("A5").Value = "jiji3"
$ws.Range("B5").Value = 34
$ws.Range("D5").Value = 25
$ws.Range("G5").Value = "牛逼3"
$ws.Range("H5").Value = "耐药"
$ws.Range("I5").Value = "耐药"
$ws.Range("J5").Value = "鸡巴"

# Date cells: write the raw serial number, then copy the number format from
# an existing date cell so the same style slot is reused instead of a new
# one being minted.
$ws.Range("C5").Value = 42970
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("E5").Value = 41631
$ws.Range("E2").Copy()
$ws.Range("E5").PasteSpecial(-4122)

# F5 is stored as a number but formatted with the text number format (same
# style already used elsewhere in the sheet, e.g. B2/B4).
$ws.Range("F5").Value = 5
$ws.Range("B2").Copy()
$ws.Range("F5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("F5").Select() | Out-Null
